# SectorGroup.xlsx edit
#
# The "codeforiati:group-code" column (originally column G, last column)
# is moved to become column D (right after the "status" column), pushing
# the former D/E/F columns (category-name, category-code, group-name)
# one position to the right, into E/F/G.
#
# In other words, for every row the 4-column block D:G is cyclically
# rotated to the right by one column:
#   new D = old G   (codeforiati:group-code)
#   new E = old D   (codeforiati:category-name)
#   new F = old E   (codeforiati:category-code)
#   new G = old F   (codeforiati:group-name)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("D1").Value = "codeforiati:group-code"
$ws.Range("E1").Value = "codeforiati:category-name"
$ws.Range("F1").Value = "codeforiati:category-code"
$ws.Range("G1").Value = "codeforiati:group-name"

# --- Data rows -----------------------------------------------------------
# Determine the last used row on the sheet.
$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

# Use far-away helper columns (Z:AC) as scratch space so the existing
# string cell type (t="s") is preserved when values are moved around
# (a plain .Value = "110" assignment would otherwise be re-interpreted
# by Excel as a number, losing the original shared-string typing).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("D$r`:G$r").Copy($ws.Range("Z$r"))
    $ws.Range("AC$r").Copy($ws.Range("D$r"))
    $ws.Range("Z$r").Copy($ws.Range("E$r"))
    $ws.Range("AA$r").Copy($ws.Range("F$r"))
    $ws.Range("AB$r").Copy($ws.Range("G$r"))
}

# Clean up the scratch area.
$ws.Range("Z2:AC" + $lastRow).Clear()
